$wb = $excel.ActiveWorkbook

# Rename the "Include from LOINC" sheet to "Include #0"
$includeSheet = $wb.Worksheets.Item("Include from LOINC")
$includeSheet.Name = "Include #0"

# Update metadata values on the "Metadata" sheet
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "4.0.2"
$ws.Range("B6").Value = "active"
$ws.Range("B8").Value = "2024-11-22T13:34:24+00:00"
